# Auto-generated edit script: applies scheduled market-data refresh to Sheets_Sagittarius_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for specific rows across sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3237
$ws.Range("I2").Value = 1279
$ws.Range("K2").Value = 1279
$ws.Range("M2").Value = -1166
$ws.Range("H39").Value = 151.14285
$ws.Range("I39").Value = 150.21053
$ws.Range("J39").Value = 160
$ws.Range("K39").Value = 450.63159
$ws.Range("L39").Value = 480
$ws.Range("M39").Value = -154.63159
$ws.Range("N39").Value = -1072
$ws.Range("H40").Value = 2225.3
$ws.Range("I40").Value = 1992.9231
$ws.Range("J40").Value = 2656.8572
$ws.Range("K40").Value = 1992.9231
$ws.Range("L40").Value = 2656.8572
$ws.Range("M40").Value = -1817.9231
$ws.Range("N40").Value = -3006.8572
$ws.Range("H80").Value = 58077.473
$ws.Range("I80").Value = 2775
$ws.Range("J80").Value = 90337.25
$ws.Range("K80").Value = 8325
$ws.Range("L80").Value = 271011.75
$ws.Range("M80").Value = -7327
$ws.Range("N80").Value = -273007.75
$ws.Range("H83").Value = 58077.473
$ws.Range("I83").Value = 2775
$ws.Range("J83").Value = 90337.25
$ws.Range("K83").Value = 24975
$ws.Range("L83").Value = 813035.25
$ws.Range("M83").Value = -19983
$ws.Range("N83").Value = -823019.25
$ws.Range("H86").Value = 3169.3076
$ws.Range("I86").Value = 1017.3333
$ws.Range("J86").Value = 5013.857
$ws.Range("K86").Value = 1017.3333
$ws.Range("L86").Value = 5013.857
$ws.Range("M86").Value = 105.6667
$ws.Range("N86").Value = -7259.857
$ws.Range("H89").Value = 3169.3076
$ws.Range("I89").Value = 1017.3333
$ws.Range("J89").Value = 5013.857
$ws.Range("K89").Value = 5086.6665
$ws.Range("L89").Value = 25069.285
$ws.Range("M89").Value = 529.3334999999997
$ws.Range("N89").Value = -36301.285
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()  # was -914
$ws.Range("H107").Value = 495.5625
$ws.Range("I107").Value = 512.5
$ws.Range("J107").Value = 444.75
$ws.Range("K107").Value = 512.5
$ws.Range("L107").Value = 444.75
$ws.Range("M107").Value = 1407.5
$ws.Range("N107").Value = -4284.75
$ws.Range("H111").Value = 1317.125
$ws.Range("I111").Value = 1089.6666
$ws.Range("J111").Value = 1999.5
$ws.Range("K111").Value = 3268.9998
$ws.Range("L111").Value = 5998.5
$ws.Range("M111").Value = -201.9998000000001
$ws.Range("N111").Value = -12132.5
$ws.Range("H125").Value = 150002220
$ws.Range("J125").Value = 166669000
$ws.Range("L125").Value = 1500021000
$ws.Range("N125").Value = -1500025920
$ws.Range("H137").Value = 3316.5
$ws.Range("I137").Value = 3125
$ws.Range("J137").Value = 3699.5
$ws.Range("K137").Value = 9375
$ws.Range("L137").Value = 11098.5
$ws.Range("M137").Value = -6825
$ws.Range("N137").Value = -16198.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9078.444
$ws.Range("I32").Value = 7004.8076
$ws.Range("K32").Value = 7004.8076
$ws.Range("M32").Value = -6717.8076

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2285.625
$ws.Range("I16").Value = 2276.2307
$ws.Range("J16").Value = 2326.3333
$ws.Range("K16").Value = 2276.2307
$ws.Range("L16").Value = 2326.3333
$ws.Range("M16").Value = -1989.2307
$ws.Range("N16").Value = -2900.3333
$ws.Range("H113").Value = 2285.625
$ws.Range("I113").Value = 2276.2307
$ws.Range("J113").Value = 2326.3333
$ws.Range("K113").Value = 2276.2307
$ws.Range("L113").Value = 2326.3333
$ws.Range("M113").Value = -106.2307000000001
$ws.Range("N113").Value = -6666.3333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 110613.5
$ws.Range("I2").Value = 100033.63
$ws.Range("J2").Value = 123544.445
$ws.Range("K2").Value = 600201.78
$ws.Range("L2").Value = 741266.67
$ws.Range("M2").Value = -600088.78
$ws.Range("N2").Value = -741492.67
$ws.Range("H17").Value = 649
$ws.Range("I17").Value = 648.5
$ws.Range("K17").Value = 1945.5
$ws.Range("M17").Value = -1776.5
$ws.Range("H34").Value = 2311.5715
$ws.Range("J34").Value = 2311.5715
$ws.Range("L34").Value = 6934.7145
$ws.Range("N34").Value = -7102.7145
$ws.Range("H39").Value = 8535.637000000001
$ws.Range("J39").Value = 8509.200000000001
$ws.Range("L39").Value = 25527.6
$ws.Range("N39").Value = -26115.6
$ws.Range("H55").Value = 1856.8182
$ws.Range("I55").Value = 2175
$ws.Range("J55").Value = 1675
$ws.Range("K55").Value = 6525
$ws.Range("L55").Value = 5025
$ws.Range("M55").Value = -6348
$ws.Range("N55").Value = -5379
$ws.Range("H80").Value = 3853.3333
$ws.Range("I80").Value = 3853.3333
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 11559.9999
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -10623.9999
$ws.Range("N80").ClearContents()  # was -7272
$ws.Range("H83").Value = 3853.3333
$ws.Range("I83").Value = 3853.3333
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 34679.9997
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -29999.9997
$ws.Range("N83").ClearContents()  # was -25560
$ws.Range("H138").Value = 3442.5715
$ws.Range("I138").Value = 3016.3333
$ws.Range("K138").Value = 9048.999899999999
$ws.Range("M138").Value = -3908.999899999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 545.25
$ws.Range("J22").Value = 275
$ws.Range("L22").Value = 275
$ws.Range("N22").Value = -1333
$ws.Range("H96").Value = 17605.5
$ws.Range("I96").Value = 15237
$ws.Range("J96").Value = 19974
$ws.Range("K96").Value = 15237
$ws.Range("L96").Value = 19974
$ws.Range("M96").Value = -12491
$ws.Range("N96").Value = -25466

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8357.727999999999
$ws.Range("I7").Value = 7576.4
$ws.Range("K7").Value = 7576.4
$ws.Range("M7").Value = -7464.4
$ws.Range("H46").Value = 40001.89
$ws.Range("I46").Value = 73860.78999999999
$ws.Range("J46").Value = 3538.4614
$ws.Range("K46").Value = 73860.78999999999
$ws.Range("L46").Value = 3538.4614
$ws.Range("M46").Value = -73672.78999999999
$ws.Range("N46").Value = -3914.4614
$ws.Range("H82").Value = 1649.3572
$ws.Range("J82").Value = 2574.75
$ws.Range("L82").Value = 2574.75
$ws.Range("N82").Value = -3296.75
$ws.Range("H85").Value = 1649.3572
$ws.Range("J85").Value = 2574.75
$ws.Range("L85").Value = 2574.75
$ws.Range("N85").Value = -5070.75
$ws.Range("H93").Value = 1075.4
$ws.Range("J93").Value = 1433.3334
$ws.Range("L93").Value = 1433.3334
$ws.Range("N93").Value = -3929.3334
$ws.Range("H126").Value = 8357.727999999999
$ws.Range("I126").Value = 7576.4
$ws.Range("K126").Value = 22729.2
$ws.Range("M126").Value = -20259.2
$ws.Range("H136").Value = 7457.778
$ws.Range("I136").Value = 8580.454
$ws.Range("K136").Value = 25741.362
$ws.Range("M136").Value = -23191.362

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 54332.668
$ws.Range("J92").Value = 54332.668
$ws.Range("L92").Value = 54332.668
$ws.Range("N92").Value = -59324.668
$ws.Range("H107").Value = 542.8
$ws.Range("I107").Value = 509.3
$ws.Range("J107").Value = 609.8
$ws.Range("K107").Value = 1527.9
$ws.Range("L107").Value = 1829.4
$ws.Range("M107").Value = 392.0999999999999
$ws.Range("N107").Value = -5669.4
$ws.Range("H126").Value = 5788.222
$ws.Range("I126").Value = 4427
$ws.Range("J126").Value = 7927.2856
$ws.Range("K126").Value = 13281
$ws.Range("L126").Value = 23781.8568
$ws.Range("M126").Value = -10811
$ws.Range("N126").Value = -28721.8568
